# Weekly update: insert this week's record at the top of the data block
# (row 9), pushing the existing rows 9-18 down to 10-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9; rows 9-18 shift down to 10-19.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with this week's data.
$ws.Cells.Item(9, 1).Value  = 7
$ws.Cells.Item(9, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value  = "Ñuble"
$ws.Cells.Item(9, 4).Value  = 44952
$ws.Cells.Item(9, 5).Value  = 16
$ws.Cells.Item(9, 6).Value  = "Fruta"
$ws.Cells.Item(9, 7).Value  = 100101
$ws.Cells.Item(9, 8).Value  = "Berries"
$ws.Cells.Item(9, 9).Value  = 100101001
$ws.Cells.Item(9, 10).Value = "Arándano (blue)"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 3000
$ws.Cells.Item(9, 15).Value = 3000
$ws.Cells.Item(9, 16).Value = 3000
$ws.Cells.Item(9, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 19).Value = 1500
$ws.Cells.Item(9, 20).Value = 2
